# Scheduled-runner price refresh: update currentAveragePrice* / LevePrice* / LeveProfit*
# columns (H:N) across the eight job sheets. Values come from an external price-API
# refresh; no formulas are involved, so each touched cell is written as a literal.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 311.75
$ws.Range("I41").Value = 332.33334
$ws.Range("J41").Value = 250
$ws.Range("K41").Value = 332.33334
$ws.Range("L41").Value = 250
$ws.Range("M41").Value = 107.66666
$ws.Range("N41").Value = -1130
$ws.Range("H53").Value = 623.73914
$ws.Range("I53").Value = 539.2308
$ws.Range("J53").Value = 733.6
$ws.Range("K53").Value = 539.2308
$ws.Range("L53").Value = 733.6
$ws.Range("M53").Value = 97.76919999999996
$ws.Range("N53").Value = -2007.6
$ws.Range("H55").Value = 176.625
$ws.Range("I55").Value = 133.4
$ws.Range("J55").Value = 248.66667
$ws.Range("K55").Value = 133.4
$ws.Range("L55").Value = 248.66667
$ws.Range("M55").Value = 80.59999999999999
$ws.Range("N55").Value = -676.6666700000001
$ws.Range("H64").Value = 12599.2
$ws.Range("I64").Value = 14330
$ws.Range("J64").Value = 10003
$ws.Range("K64").Value = 14330
$ws.Range("L64").Value = 10003
$ws.Range("M64").Value = -14082
$ws.Range("N64").Value = -10499
$ws.Range("H67").Value = 12599.2
$ws.Range("I67").Value = 14330
$ws.Range("J67").Value = 10003
$ws.Range("K67").Value = 14330
$ws.Range("L67").Value = 10003
$ws.Range("M67").Value = -13472
$ws.Range("N67").Value = -11719
$ws.Range("H92").Value = 1954.4
$ws.Range("I92").Value = 1838.3334
$ws.Range("K92").Value = 1838.3334
$ws.Range("M92").Value = -590.3334
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H112").Value = 10206445
$ws.Range("J112").Value = 10871826
$ws.Range("L112").Value = 32615478
$ws.Range("N112").Value = -32617694
$ws.Range("H132").Value = 3829.9644
$ws.Range("I132").Value = 3808.8518
$ws.Range("K132").Value = 11426.5554
$ws.Range("M132").Value = -8896.555399999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 34486780
$ws.Range("I61").Value = 47621856
$ws.Range("K61").Value = 47621856
$ws.Range("M61").Value = -47621644
$ws.Range("H74").Value = 22223376
$ws.Range("I74").Value = 25001048
$ws.Range("K74").Value = 25001048
$ws.Range("M74").Value = -25000174
$ws.Range("H77").Value = 22223376
$ws.Range("I77").Value = 25001048
$ws.Range("K77").Value = 125005240
$ws.Range("M77").Value = -125000872
$ws.Range("H110").Value = 1266.8695
$ws.Range("I110").Value = 558.7
$ws.Range("K110").Value = 558.7
$ws.Range("M110").Value = 1486.3
$ws.Range("H132").Value = 24456126
$ws.Range("I132").Value = 8275.161
$ws.Range("J132").Value = 100244460
$ws.Range("K132").Value = 24825.483
$ws.Range("L132").Value = 300733380
$ws.Range("M132").Value = -22295.483
$ws.Range("N132").Value = -300738440
$ws.Range("H136").Value = 34486780
$ws.Range("I136").Value = 47621856
$ws.Range("K136").Value = 142865568
$ws.Range("M136").Value = -142863018

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1337.4286
$ws.Range("I64").Value = 1722.2222
$ws.Range("J64").Value = 1048.8334
$ws.Range("K64").Value = 1722.2222
$ws.Range("L64").Value = 1048.8334
$ws.Range("M64").Value = -1497.2222
$ws.Range("N64").Value = -1498.8334
$ws.Range("H67").Value = 1337.4286
$ws.Range("I67").Value = 1722.2222
$ws.Range("J67").Value = 1048.8334
$ws.Range("K67").Value = 1722.2222
$ws.Range("L67").Value = 1048.8334
$ws.Range("M67").Value = -942.2221999999999
$ws.Range("N67").Value = -2608.8334
$ws.Range("H99").Value = 3462.8286
$ws.Range("I99").Value = 2925.5186
$ws.Range("J99").Value = 5276.25
$ws.Range("K99").Value = 2925.5186
$ws.Range("L99").Value = 5276.25
$ws.Range("M99").Value = -1427.5186
$ws.Range("N99").Value = -8272.25
$ws.Range("H134").Value = 1982.3
$ws.Range("I134").Value = 1466.8182
$ws.Range("J134").Value = 4412.4287
$ws.Range("K134").Value = 4400.4546
$ws.Range("L134").Value = 13237.2861
$ws.Range("M134").Value = -1865.4546
$ws.Range("N134").Value = -18307.2861

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 15877615
$ws.Range("I31").Value = 3329.7273
$ws.Range("K31").Value = 3329.7273
$ws.Range("M31").Value = -3034.7273
$ws.Range("H34").Value = 15877615
$ws.Range("I34").Value = 3329.7273
$ws.Range("K34").Value = 3329.7273
$ws.Range("M34").Value = -3127.7273
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("N45").ClearContents()
$ws.Range("H53").Value = 219430.75
$ws.Range("J53").Value = 274637.5
$ws.Range("L53").Value = 274637.5
$ws.Range("N53").Value = -275851.5
$ws.Range("H58").Value = 2190.5278
$ws.Range("I58").Value = 1871.8214
$ws.Range("K58").Value = 1871.8214
$ws.Range("M58").Value = -1668.8214
$ws.Range("H132").Value = 43411.625
$ws.Range("I132").Value = 48318.094
$ws.Range("K132").Value = 144954.282
$ws.Range("M132").Value = -142424.282
$ws.Range("H134").Value = 2982.72
$ws.Range("I134").Value = 2790.1702
$ws.Range("K134").Value = 8370.5106
$ws.Range("M134").Value = -5835.5106
$ws.Range("H136").Value = 2190.5278
$ws.Range("I136").Value = 1871.8214
$ws.Range("K136").Value = 5615.4642
$ws.Range("M136").Value = -3065.4642

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2121.76
$ws.Range("I5").Value = 1079.8
$ws.Range("J5").Value = 2382.25
$ws.Range("K5").Value = 3239.4
$ws.Range("L5").Value = 7146.75
$ws.Range("M5").Value = -3127.4
$ws.Range("N5").Value = -7370.75
$ws.Range("H8").Value = 349
$ws.Range("I8").Value = 349
$ws.Range("K8").Value = 1047
$ws.Range("M8").Value = -908
$ws.Range("H25").Value = 3999.3333
$ws.Range("J25").Value = 3999.3333
$ws.Range("L25").Value = 11997.9999
$ws.Range("N25").Value = -12335.9999
$ws.Range("H30").Value = 3999.3333
$ws.Range("J30").Value = 3999.3333
$ws.Range("L30").Value = 11997.9999
$ws.Range("N30").Value = -12201.9999
$ws.Range("H130").Value = 2691.5
$ws.Range("I130").Value = 1500
$ws.Range("J130").Value = 3088.6667
$ws.Range("K130").Value = 4500
$ws.Range("L130").Value = 9266.000100000001
$ws.Range("M130").Value = 520
$ws.Range("N130").Value = -19306.0001
$ws.Range("H132").Value = 3925493.8
$ws.Range("I132").Value = 2158.4
$ws.Range("K132").Value = 19425.6
$ws.Range("M132").Value = -16895.6
$ws.Range("H135").Value = 2121.76
$ws.Range("I135").Value = 1079.8
$ws.Range("J135").Value = 2382.25
$ws.Range("K135").Value = 9718.199999999999
$ws.Range("L135").Value = 21440.25
$ws.Range("M135").Value = -7183.199999999999
$ws.Range("N135").Value = -26510.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
$ws.Range("H132").Value = 1541.8611
$ws.Range("I132").Value = 1573.7241
$ws.Range("K132").Value = 4721.1723
$ws.Range("M132").Value = -2191.1723

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4264.357
$ws.Range("I16").Value = 4264.357
$ws.Range("K16").Value = 4264.357
$ws.Range("M16").Value = -4094.357
$ws.Range("H46").Value = 2076.2307
$ws.Range("I46").Value = 776.7778
$ws.Range("K46").Value = 776.7778
$ws.Range("M46").Value = -588.7778
$ws.Range("H55").Value = 754.2941
$ws.Range("I55").Value = 465.3
$ws.Range("J55").Value = 1167.1428
$ws.Range("K55").Value = 465.3
$ws.Range("L55").Value = 1167.1428
$ws.Range("M55").Value = -292.3
$ws.Range("N55").Value = -1513.1428
$ws.Range("H122").Value = 6911.375
$ws.Range("J122").Value = 8418.200000000001
$ws.Range("L122").Value = 25254.6
$ws.Range("N122").Value = -30154.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 62001.5
$ws.Range("J70").Value = 62001.5
$ws.Range("L70").Value = 62001.5
$ws.Range("N70").Value = -62631.5
$ws.Range("H73").Value = 62001.5
$ws.Range("J73").Value = 62001.5
$ws.Range("L73").Value = 62001.5
$ws.Range("N73").Value = -64185.5
$ws.Range("I122").Value = 125125624
$ws.Range("J122").Value = 3140.7778
$ws.Range("K122").Value = 375376872
$ws.Range("L122").Value = 9422.3334
$ws.Range("M122").Value = -375374422
$ws.Range("N122").Value = -14322.3334
